$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.326.53"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "'3.778.34"
$ws.Range("E3").Value = "  +2.39%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'410.38"
$ws.Range("E5").Value = "  -2.16%  "
$ws.Range("D6").Value = "'132.85"
$ws.Range("D7").Value = "'3.769.93"
$ws.Range("E7").Value = "  +2.48%  "
$ws.Range("D8").Value = "'0.615"
$ws.Range("E8").Value = "  -4.51%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.16%  "
$ws.Range("D10").Value = "'0.736"
$ws.Range("E10").Value = "  -3.04%  "
$ws.Range("D11").Value = "'0.167"
$ws.Range("E11").Value = "  -7.82%  "
$ws.Range("D12").Value = "'0.0000364"
$ws.Range("E12").Value = "  -7.12%  "
$ws.Range("D13").Value = "'41.07"
$ws.Range("E13").Value = "  -4.11%  "
$ws.Range("D14").Value = "'10.16"
$ws.Range("E14").Value = "  -4.05%  "
$ws.Range("D15").Value = "'4.353.62"
$ws.Range("E15").Value = "  +1.81%  "
$ws.Range("D16").Value = "'14.76"
$ws.Range("E16").Value = "  +12.57%  "
$ws.Range("E17").Value = "  -1.31%  "
$ws.Range("D18").Value = "'3.774.96"
$ws.Range("E18").Value = "  +2.29%  "
$ws.Range("E19").Value = "  -5.03%  "
$ws.Range("D20").Value = "'66.384.74"
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("E21").Value = "  -4.84%  "
$ws.Range("D22").Value = "'413.29"
$ws.Range("E22").Value = "  -6.73%  "
$ws.Range("D23").Value = "'14.56"
$ws.Range("E23").Value = "  -10.63%  "
$ws.Range("D24").Value = "'85.64"
$ws.Range("E24").Value = "  -4.41%  "
$ws.Range("E25").Value = "  -1.69%  "
$ws.Range("D26").Value = "'5.72"
$ws.Range("E26").Value = "  +14.51%  "
$ws.Range("D27").Value = "'36.12"
$ws.Range("E27").Value = "  -3.72%  "
$ws.Range("E28").Value = "  -5.24%  "
$ws.Range("D29").Value = "'9.44"
$ws.Range("E29").Value = "  -8.11%  "
$ws.Range("D30").Value = "'707.27"
$ws.Range("E30").Value = "  +8.67%  "
$ws.Range("E31").Value = "  -0.44%  "
$ws.Range("D32").Value = "'12.42"
$ws.Range("E32").Value = "  -2.08%  "
$ws.Range("D33").Value = "'2.75"
$ws.Range("E33").Value = "  -0.60%  "
$ws.Range("D34").Value = "'7.38"
$ws.Range("E34").Value = "  +1.59%  "
$ws.Range("D35").Value = "'0.155"
$ws.Range("E35").Value = "  -5.60%  "
$ws.Range("D36").Value = "'39.23"
$ws.Range("E36").Value = "  -4.96%  "
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").Value = "'55.24"
$ws.Range("E38").Value = "  -3.50%  "
$ws.Range("D39").Value = "'0.0₃0743"
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("D40").Value = "'0.0462"
$ws.Range("E40").Value = "  -6.61%  "
$ws.Range("D41").Value = "'2.86"
$ws.Range("E41").Value = "  -10.41%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").Value = "'0.137"
$ws.Range("E43").Value = "  -8.30%  "
$ws.Range("D44").Value = "'27.37"
$ws.Range("E44").Value = "  -5.82%  "
$ws.Range("E45").Value = "  +19.15%  "
$ws.Range("D46").Value = "'146.10"
$ws.Range("E46").Value = "  -1.32%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "'2.08"
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("B48").Value = "LidoDAOToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D48").Value = "'3.28"
$ws.Range("E48").Value = "  -3.67%  "
$ws.Range("D49").Value = "'4.29"
$ws.Range("E49").Value = "  -1.89%  "
$ws.Range("E50").Value = "  +0.48%  "
$ws.Range("E51").Value = "  -3.35%  "
